$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The old last row of the "Diario de bordo" table (date 09/04) had a
#    "_GoBack" bookmark at the end of its second cell's last paragraph.
#    That bookmark is moving to the new row we are about to append, so we
#    strip it from here first (while paragraph offsets are still pristine).
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$targetPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Iniciado mecânicas*") {
        $targetPara = $p
        break
    }
}

$bookmarkParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>Iniciado mecânicas de movimentação no Script Player.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> Iniciado máquina de estados e criado parâmetros de controle para </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Idle</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> e </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>run</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t>.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

if ($targetPara -ne $null) {
    $targetPara.Range.InsertXML($bookmarkParaXml)
}

# ---------------------------------------------------------------------------
# 2) Append the new "10/04" row to the first table.
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$newRow = $t1.Rows.Add()

# --- Cell 1: date -----------------------------------------------------
$newRow.Cells.Item(1).Range.Text = "10/04"

# --- Cell 2: the day's log, with the _GoBack bookmark at the very end -
$cell2 = $newRow.Cells.Item(2)
$cell2Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">Criado </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Colliders</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> do player e Tiles, Implementado </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Level</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> inicial para teste, criado mecânica inic</w:t></w:r>
            <w:r><w:t>i</w:t></w:r>
            <w:r><w:t xml:space="preserve">al de pulo, configurado </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Tilemap</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> e escalonamento inicial de personagem/tile/</w:t></w:r>
            <w:r><w:t>câmera. Esboço do primeiro inimigo</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$cell2.Range.InsertXML($cell2Xml)
# The cell started with one empty placeholder paragraph; InsertXML appended
# the real content as a second paragraph, so drop the leading empty one.
$cell2.Range.Paragraphs.Item(1).Range.Delete()

# --- Cell 3: "Nao" (misspelled by the author, hence the proofErr wrap) -
$cell3 = $newRow.Cells.Item(3)
$cell3Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Nao</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$cell3.Range.InsertXML($cell3Xml)
$cell3.Range.Paragraphs.Item(1).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Second (Gantt-style) table: recolor a handful of status cells.
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

function Set-CellShade($table, $rowIndex, $cellIndex, $fill, $themeFill) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item($cellIndex)
    $cell.Shading.BackgroundPatternColor = $fill
    if ($themeFill -ne $null) {
        $cell.Shading.Texture = $themeFill
    }
}

# Row 2 = "Selecionar/desenhar a arte dos personagens"
$t2.Rows.Item(2).Cells.Item(2).Shading.BackgroundPatternColor = 0x47AD70   ; # 70AD47 (BGR order)
$t2.Rows.Item(2).Cells.Item(3).Shading.BackgroundPatternColor = 0x00FFFF   ; # FFFF00 (BGR order)

# Row 3 = "Selecionar/desenhar a arte dos cenários"
$t2.Rows.Item(3).Cells.Item(2).Shading.BackgroundPatternColor = 0x47AD70
$t2.Rows.Item(3).Cells.Item(3).Shading.BackgroundPatternColor = 0x00FFFF

# Row 4 = "Desenvolver o sistema de controle do jogador"
$t2.Rows.Item(4).Cells.Item(4).Shading.BackgroundPatternColor = 0x47AD70

# Row 5 = "Desenvolver Tela e sistema de Colisão"
$t2.Rows.Item(5).Cells.Item(4).Shading.BackgroundPatternColor = 0x47AD70

Write-Host "done"
